$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Reference style to apply so the GUID cells stay text (preserving leading zeros)
# while keeping the same visual style as the rest of the directory listing.
$refStyleA = $ws.Range("A87").Style
$refStyleB = $ws.Range("B87").Style

# Row 88: new GUID entry 000087
$ws.Range("A88").NumberFormat = "@"
$ws.Range("A88").Value = "000087"
$ws.Range("A88").Style = $refStyleA
$ws.Range("B88").Value = "Details: Baton Tip Pose Transformation. IMU CJMCU-20948 Data Reading, Fused with imufilter, transformed with BatonTip_Transformation. Script used: BatonTipPoseVisualisation.  Dataset used: IMU data: IMU_Orientation_Reading_08_03_23. Transformed Baton tip data: BatonTipPose_08_03_23.. File Location: Visualisations/IMU_TransformedBatonTipPose. Date Generated: 08-Mar-2023 13:46:35"
$ws.Range("B88").Style = $refStyleB

# Row 89: new GUID entry 000088
$ws.Range("A89").NumberFormat = "@"
$ws.Range("A89").Value = "000088"
$ws.Range("A89").Style = $refStyleA
$ws.Range("B89").Value = "Details: Baton Tip Pose Transformation - IMU CJMCU-20948 Data Reading, Fused with imufilter, transformed with BatonTip_Transformation. Script used: BatonTipPoseVisualisation.  Dataset used: IMU data: IMU_Orientation_Reading_08_03_23. Transformed Baton tip data: BatonTipPose_08_03_23.. File Location: Visualisations/IMU_TransformedBatonTipPose. Date Generated: 08-Mar-2023 13:46:56"
$ws.Range("B89").Style = $refStyleB

# Row 90: new GUID entry 000089
$ws.Range("A90").NumberFormat = "@"
$ws.Range("A90").Value = "000089"
$ws.Range("A90").Style = $refStyleA
$ws.Range("B90").Value = "Details: Baton Tip Pose Transformation - IMU CJMCU-20948 Data Reading - Fused with imufilter - transformed with BatonTip_Transformation. Script used: BatonTipPoseVisualisation.  Dataset used: IMU data: IMU_Orientation_Reading_08_03_23. Transformed Baton tip data: BatonTipPose_08_03_23.. File Location: Visualisations/IMU_TransformedBatonTipPose. Date Generated: 08-Mar-2023 13:48:33"
$ws.Range("B90").Style = $refStyleB

# Column B widened to fit the new, longer content
$ws.Columns.Item(2).ColumnWidth = 349.5924479166667
